# Added fixed headers for vertical scrolling
# (Underlying data change: extend the Inventory table with 20 more rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# Fill in the new rows 12..31 -> Column A = 11..30, Column B = "J"
for ($i = 12; $i -le 31; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 2).Value = "J"
}

# Resize the table ("Library") so it covers the new data range
$table = $ws.ListObjects.Item("Library")
$table.Resize($ws.Range("A1:K31"))

# Update the active selection to reflect scrolling down to view new rows
$ws.Activate()
$ws.Range("F21").Select()
